$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Sending cluster" / "Target cluster" data now includes a third cluster (ECs),
# expanding the table from 6 data rows (2 clusters x 3 combos) to 9 rows (3 clusters x 3 combos).
$data = New-Object 'object[,]' 9,20
$data[0,0] = "ECs"
$data[0,1] = "Bmp2"
$data[0,2] = "Eng"
$data[0,3] = "ECs"
$data[0,4] = 2
$data[0,5] = 0.6666666666666666
$data[0,6] = 1.445484
$data[0,7] = 4.336452
$data[0,8] = 0.1286708197254238
$data[0,9] = 0.1286708197254238
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 170.93328
$data[0,13] = 512.79984
$data[0,14] = 0.7687311215213114
$data[0,15] = 0.7687311215213115
$data[0,16] = 247.0813213075199
$data[0,17] = 2223.73189176768
$data[0,18] = 0.09891326355459153
$data[0,19] = 0.09891326355459155
$data[1,0] = "ECs"
$data[1,1] = "Bmp2"
$data[1,2] = "Eng"
$data[1,3] = "FAPs"
$data[1,4] = 2
$data[1,5] = 0.6666666666666666
$data[1,6] = 1.445484
$data[1,7] = 4.336452
$data[1,8] = 0.1286708197254238
$data[1,9] = 0.1286708197254238
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 40.31217066666667
$data[1,13] = 120.936512
$data[1,14] = 0.1812942463137967
$data[1,15] = 0.1812942463137967
$data[1,16] = 58.27059770393599
$data[1,17] = 524.4353793354239
$data[1,18] = 0.02332727928469912
$data[1,19] = 0.02332727928469912
$data[2,0] = "ECs"
$data[2,1] = "Bmp2"
$data[2,2] = "Eng"
$data[2,3] = "sCs"
$data[2,4] = 2
$data[2,5] = 0.6666666666666666
$data[2,6] = 1.445484
$data[2,7] = 4.336452
$data[2,8] = 0.1286708197254238
$data[2,9] = 0.1286708197254238
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 11.112244
$data[2,13] = 33.336732
$data[2,14] = 0.04997463216489184
$data[2,15] = 0.04997463216489184
$data[2,16] = 16.062570906096
$data[2,17] = 144.563138154864
$data[2,18] = 0.006430276886133165
$data[2,19] = 0.006430276886133165
$data[3,0] = "FAPs"
$data[3,1] = "Bmp2"
$data[3,2] = "Eng"
$data[3,3] = "ECs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 6.292313
$data[3,7] = 18.876939
$data[3,8] = 0.5601148623429528
$data[3,9] = 0.5601148623429528
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 170.93328
$data[3,13] = 512.79984
$data[3,14] = 0.7687311215213114
$data[3,15] = 0.7687311215213115
$data[3,16] = 1075.56569987664
$data[3,17] = 9680.09129888976
$data[3,18] = 0.4305777263096531
$data[3,19] = 0.4305777263096531
$data[4,0] = "FAPs"
$data[4,1] = "Bmp2"
$data[4,2] = "Eng"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 6.292313
$data[4,7] = 18.876939
$data[4,8] = 0.5601148623429528
$data[4,9] = 0.5601148623429528
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 40.31217066666667
$data[4,13] = 120.936512
$data[4,14] = 0.1812942463137967
$data[4,15] = 0.1812942463137967
$data[4,16] = 253.6567955440853
$data[4,17] = 2282.911159896768
$data[4,18] = 0.1015456018176216
$data[4,19] = 0.1015456018176216
$data[5,0] = "FAPs"
$data[5,1] = "Bmp2"
$data[5,2] = "Eng"
$data[5,3] = "sCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 6.292313
$data[5,7] = 18.876939
$data[5,8] = 0.5601148623429528
$data[5,9] = 0.5601148623429528
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 11.112244
$data[5,13] = 33.336732
$data[5,14] = 0.04997463216489184
$data[5,15] = 0.04997463216489184
$data[5,16] = 69.92171738037199
$data[5,17] = 629.295456423348
$data[5,18] = 0.0279915342156781
$data[5,19] = 0.0279915342156781
$data[6,0] = "sCs"
$data[6,1] = "Bmp2"
$data[6,2] = "Eng"
$data[6,3] = "ECs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 3.496172
$data[6,7] = 10.488516
$data[6,8] = 0.3112143179316233
$data[6,9] = 0.3112143179316232
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 170.93328
$data[6,13] = 512.79984
$data[6,14] = 0.7687311215213114
$data[6,15] = 0.7687311215213115
$data[6,16] = 597.6121474041599
$data[6,17] = 5378.509326637441
$data[6,18] = 0.2392401316570668
$data[6,19] = 0.2392401316570667
$data[7,0] = "sCs"
$data[7,1] = "Bmp2"
$data[7,2] = "Eng"
$data[7,3] = "FAPs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 3.496172
$data[7,7] = 10.488516
$data[7,8] = 0.3112143179316233
$data[7,9] = 0.3112143179316232
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 40.31217066666667
$data[7,13] = 120.936512
$data[7,14] = 0.1812942463137967
$data[7,15] = 0.1812942463137967
$data[7,16] = 140.9382823440213
$data[7,17] = 1268.444541096192
$data[7,18] = 0.05642136521147595
$data[7,19] = 0.05642136521147594
$data[8,0] = "sCs"
$data[8,1] = "Bmp2"
$data[8,2] = "Eng"
$data[8,3] = "sCs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 3.496172
$data[8,7] = 10.488516
$data[8,8] = 0.3112143179316233
$data[8,9] = 0.3112143179316232
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 11.112244
$data[8,13] = 33.336732
$data[8,14] = 0.04997463216489184
$data[8,15] = 0.04997463216489184
$data[8,16] = 38.850316329968
$data[8,17] = 349.652846969712
$data[8,18] = 0.01555282106308058
$data[8,19] = 0.01555282106308058

$ws.Range("A2:T10").Value = $data
